$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (D, E, F) before the existing "Terms Typically Offered" column,
# shifting it to column G, and set the new header row.
$ws.Range("D1:F1").EntireColumn.Insert()

$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Row-by-row data: update Prerequisites (C) text to remove the embedded
# corequisite/concurrent/recommended clauses, and populate the three new columns
# plus the shifted Terms Typically Offered column (G) accordingly.

$ws.Range("C2").Value = "MATE majors only."
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "F"

$ws.Range("C3").Value = "MATE 110."
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "W"

$ws.Range("C4").Value = "MATE 120."
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "SP"

$ws.Range("C5").Value = "Consent of department head."
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "F,W,SP,SU"

$ws.Range("C6").Value = "CHEM 111 or CHEM 124 or CHEM 127."
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "Concurrent enrollment in MATE 215."
$ws.Range("G6").Value = "F,W,SP,SU "

$ws.Range("C7").Value = "NA"
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("G7").Value = "F,W,SP,SUPrerequisite or MATE 210."

$ws.Range("C8").Value = "MATE 210."
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "W"

$ws.Range("C9").Value = "MATE 215."
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "MATE 232."
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "F "

$ws.Range("C10").Value = "MATE 210."
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "F"

$ws.Range("C11").Value = "MATE 225."
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"
$ws.Range("G11").Value = "W"

$ws.Range("C12").Value = "Open to undergraduate students and consent of instructor."
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "TBD"

$ws.Range("C13").Value = "CHEM 125, PHYS 133, MATH 143, MATE 210 and MATE 215."
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("G13").Value = "W"

$ws.Range("C14").Value = "MATE 210, MATE 340, STAT 312."
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "MATE 350."
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "SP "

$ws.Range("C15").Value = "MATE 360."
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "W"

$ws.Range("C16").Value = "PHYS 132 and MATH 141."
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "F"

$ws.Range("C17").Value = "MATH 141 and ME 211."
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "W"

$ws.Range("C18").Value = "CHEM 124 and PHYS 133."
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "SP"

$ws.Range("C19").Value = "MATE 210 and PHYS 133."
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "NA"
$ws.Range("G19").Value = "F"

$ws.Range("C20").Value = "MATE 360."
$ws.Range("D20").Value = "CE 204."
$ws.Range("E20").Value = "MATE 310."
$ws.Range("F20").Value = "NA"
$ws.Range("G20").Value = "SP  "

$ws.Range("C21").Value = "Junior standing; completion of GE Area A with grades of C- or better; and completion of one course in GE Area B1 with a grade of C- or better."
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "F,W,SP,SU"

$ws.Range("C22").Value = "MATE 232 and MATE 235."
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "SP"

$ws.Range("C23").Value = "MATE 360 and MATE 380."
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("F23").Value = "NA"
$ws.Range("G23").Value = "F"

$ws.Range("C24").Value = "CHEM 125, PHYS 133, MATH 143, MATE 210 and MATE 215; Materials Engineering students only."
$ws.Range("D24").Value = "NA"
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "W"

$ws.Range("C25").Value = "Consent of department head."
$ws.Range("D25").Value = "NA"
$ws.Range("E25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "F,W,SP,SU"

$ws.Range("C26").Value = "MATE 210 and MATE 215."
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = "NA"
$ws.Range("F26").Value = "NA"
$ws.Range("G26").Value = "F, W, SP"

$ws.Range("C27").Value = "MATE 210."
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = "NA"
$ws.Range("F27").Value = "MATE 401."
$ws.Range("G27").Value = "F, W, SP "

$ws.Range("C28").Value = "CSC 231, ME 211, MATE 280."
$ws.Range("D28").Value = "NA"
$ws.Range("E28").Value = "NA"
$ws.Range("F28").Value = "NA"
$ws.Range("G28").Value = "F, W, SP"

$ws.Range("C29").Value = "CHEM 125, PHYS 133 and MATE 210."
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "TBD"

$ws.Range("C30").Value = "MATE 310 or CHEM 444."
$ws.Range("D30").Value = "NA"
$ws.Range("E30").Value = "NA"
$ws.Range("F30").Value = "NA"
$ws.Range("G30").Value = "W"

$ws.Range("C31").Value = "MATE 210."
$ws.Range("D31").Value = "NA"
$ws.Range("E31").Value = "NA"
$ws.Range("F31").Value = "NA"
$ws.Range("G31").Value = "W"

$ws.Range("C32").Value = "CHEM 125 or CHEM 128, MATE 210, MATE 215."
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = "NA"
$ws.Range("F32").Value = "NA"
$ws.Range("G32").Value = "F"

$ws.Range("C33").Value = "BMED 212 or MATE 210."
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = "NA"
$ws.Range("F33").Value = "NA"
$ws.Range("G33").Value = "W"

$ws.Range("C34").Value = "NA"
$ws.Range("D34").Value = "BMED 434/EE 423/MATE 430."
$ws.Range("E34").Value = "NA"
$ws.Range("F34").Value = "NA"
$ws.Range("G34").Value = "W"

$ws.Range("C35").Value = "MATE 210."
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = "NA"
$ws.Range("F35").Value = "NA"
$ws.Range("G35").Value = "SP"

$ws.Range("C36").Value = "MATE 210."
$ws.Range("D36").Value = "MATE 440."
$ws.Range("E36").Value = "NA"
$ws.Range("F36").Value = "NA"
$ws.Range("G36").Value = "SP "

$ws.Range("C37").Value = "CHEM 125 or CHEM 128; CHEM 351, MATE 380, or ME 302."
$ws.Range("D37").Value = "NA"
$ws.Range("E37").Value = "NA"
$ws.Range("F37").Value = "NA"
$ws.Range("G37").Value = "SP"

$ws.Range("C38").Value = "MATE 210, MATE 215, MATE 350; and Senior standing."
$ws.Range("D38").Value = "NA"
$ws.Range("E38").Value = "NA"
$ws.Range("F38").Value = "NA"
$ws.Range("G38").Value = "W"

$ws.Range("C39").Value = "MATE 210."
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = "NA"
$ws.Range("F39").Value = "NA"
$ws.Range("G39").Value = "TBD"

$ws.Range("C40").Value = "EE 112 or EE 113 or EE 201."
$ws.Range("D40").Value = "NA"
$ws.Range("E40").Value = "NA"
$ws.Range("F40").Value = "MATE 210."
$ws.Range("G40").Value = "F, W, SP "

$ws.Range("C41").Value = "MATE 350."
$ws.Range("D41").Value = "NA"
$ws.Range("E41").Value = "NA"
$ws.Range("F41").Value = "NA"
$ws.Range("G41").Value = "TBD"

$ws.Range("C42").Value = "MATE 360."
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("F42").Value = "NA"
$ws.Range("G42").Value = "F"

$ws.Range("C43").Value = "Consent of instructor."
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "NA"
$ws.Range("F43").Value = "NA"
$ws.Range("G43").Value = "TBD"

$ws.Range("C44").Value = "Consent of instructor."
$ws.Range("D44").Value = "NA"
$ws.Range("E44").Value = "NA"
$ws.Range("F44").Value = "NA"
$ws.Range("G44").Value = "TBD"

$ws.Range("C45").Value = "MATE 350."
$ws.Range("D45").Value = "NA"
$ws.Range("E45").Value = "NA"
$ws.Range("F45").Value = "NA"
$ws.Range("G45").Value = "F"

$ws.Range("C46").Value = "IME 144; senior standing; and Materials Engineering major."
$ws.Range("D46").Value = "NA"
$ws.Range("E46").Value = "NA"
$ws.Range("F46").Value = "NA"
$ws.Range("G46").Value = "F"

$ws.Range("C47").Value = "MATE 482."
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "W"

$ws.Range("C48").Value = "MATE 483."
$ws.Range("D48").Value = "NA"
$ws.Range("E48").Value = "NA"
$ws.Range("F48").Value = "NA"
$ws.Range("G48").Value = "SP"

$ws.Range("C49").Value = "MATE 210."
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = "NA"
$ws.Range("F49").Value = "NA"
$ws.Range("G49").Value = "TBD"

$ws.Range("C50").Value = "MATE 210."
$ws.Range("D50").Value = "NA"
$ws.Range("E50").Value = "NA"
$ws.Range("F50").Value = "NA"
$ws.Range("G50").Value = "W"

$ws.Range("C51").Value = "Sophomore standing and consent of instructor."
$ws.Range("D51").Value = "NA"
$ws.Range("E51").Value = "NA"
$ws.Range("F51").Value = "NA"
$ws.Range("G51").Value = "F,W,SP,SU"

$ws.Range("C52").Value = "Consent of department head, graduate advisor, or supervising faculty member."
$ws.Range("D52").Value = "NA"
$ws.Range("E52").Value = "NA"
$ws.Range("F52").Value = "NA"
$ws.Range("G52").Value = "F,W,SP,SU"

$ws.Range("C53").Value = "BIO 161, or BIO 213 and BMED/BRAE 213; MATE 210 and graduate standing or consent of instructor."
$ws.Range("D53").Value = "NA"
$ws.Range("E53").Value = "NA"
$ws.Range("F53").Value = "NA"
$ws.Range("G53").Value = "W"

$ws.Range("C54").Value = "Graduate standing."
$ws.Range("D54").Value = "NA"
$ws.Range("E54").Value = "NA"
$ws.Range("F54").Value = "NA"
$ws.Range("G54").Value = "F, W, SP"

$ws.Range("C55").Value = "Graduate standing or consent of instructor."
$ws.Range("D55").Value = "NA"
$ws.Range("E55").Value = "NA"
$ws.Range("F55").Value = "NA"
$ws.Range("G55").Value = "TBD"

$ws.Range("C56").Value = "Senior or graduate standing or consent of instructor."
$ws.Range("D56").Value = "NA"
$ws.Range("E56").Value = "NA"
$ws.Range("F56").Value = "NA"
$ws.Range("G56").Value = "TBD"

$ws.Range("C57").Value = "Graduate standing."
$ws.Range("D57").Value = "NA"
$ws.Range("E57").Value = "NA"
$ws.Range("F57").Value = "NA"
$ws.Range("G57").Value = "F,W,SP,SU"
